$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

$ws.Range("D2").Value = "26.201.64"
$ws.Range("E2").Value = "  -1.98%  "
$ws.Range("D3").Value = "1.672.44"
$ws.Range("E3").Value = "  -1.35%  "
Set-TextCell "D4" "1.006"
$ws.Range("E4").Value = "  +0.07%  "
Set-TextCell "D5" "217.31"
$ws.Range("E5").Value = "  -0.85%  "
Set-TextCell "D6" "0.5124"
$ws.Range("E6").Value = "  +0.87%  "
Set-TextCell "D7" "1.006"
$ws.Range("E7").Value = "  +0.06%  "
Set-TextCell "D8" "0.2638"
$ws.Range("E8").Value = "  +1.52%  "
Set-TextCell "D9" "0.06402"
$ws.Range("E9").Value = "  +4.28%  "
Set-TextCell "D10" "21.61"
$ws.Range("E10").Value = "  -1.27%  "
Set-TextCell "D11" "0.07403"
$ws.Range("E11").Value = "  +0.86%  "
$ws.Range("D12").Value = "1.684.48"
$ws.Range("E12").Value = "  -0.54%  "
Set-TextCell "D13" "4.522"
$ws.Range("E13").Value = "  +2.19%  "
Set-TextCell "D14" "0.5818"
$ws.Range("E14").Value = "  +1.39%  "
Set-TextCell "D15" "0.000008599"
$ws.Range("E15").Value = "  +5.51%  "
Set-TextCell "D16" "64.46"
$ws.Range("E16").Value = "  -1.59%  "
$ws.Range("D17").Value = "26.252.70"
$ws.Range("E17").Value = "  -1.94%  "
Set-TextCell "D18" "4.940"
$ws.Range("E18").Value = "  -1.42%  "
$ws.Range("E19").Value = "  +0.09%  "
Set-TextCell "D20" "10.83"
$ws.Range("E20").Value = "  +1.17%  "
Set-TextCell "D21" "190.21"
$ws.Range("E21").Value = "  +3.16%  "
Set-TextCell "D22" "6.215"
$ws.Range("E22").Value = "  -0.01%  "
$ws.Range("E23").Value = "  +0.11%  "
Set-TextCell "D24" "145.61"
$ws.Range("E24").Value = "  +0.34%  "
Set-TextCell "D25" "7.624"
$ws.Range("E25").Value = "  -0.32%  "
Set-TextCell "D26" "0.1188"
$ws.Range("E26").Value = "  +3.61%  "
$ws.Range("E27").Value = "  +2.88%  "
Set-TextCell "D28" "0.06332"
$ws.Range("E28").Value = "  +10.88%  "
Set-TextCell "D29" "1.298"
$ws.Range("E29").Value = "  -1.45%  "
Set-TextCell "D30" "1.322"
$ws.Range("E30").Value = "  -1.05%  "
Set-TextCell "D31" "3.522"
$ws.Range("E31").Value = "  +1.34%  "
Set-TextCell "D32" "3.520"
$ws.Range("E32").Value = "  +2.80%  "
Set-TextCell "D33" "1.634"
$ws.Range("E33").Value = "  -1.78%  "
Set-TextCell "D34" "1.017"
$ws.Range("E34").Value = "  +1.24%  "
Set-TextCell "D35" "0.6070"
$ws.Range("E35").Value = "  +2.82%  "
Set-TextCell "D36" "2.381"
$ws.Range("E36").Value = "  -1.35%  "
Set-TextCell "D37" "2.661"
$ws.Range("E37").Value = "  +0.66%  "
Set-TextCell "D38" "6.178"
$ws.Range("E38").Value = "  +4.41%  "
Set-TextCell "D39" "0.01606"
$ws.Range("E39").Value = "  +0.86%  "
$ws.Range("D40").Value = "1.084.32"
$ws.Range("E40").Value = "  +1.48%  "
Set-TextCell "D41" "0.8656"
$ws.Range("E41").Value = "  +1.20%  "
$ws.Range("E42").Value = "  +0.69%  "
Set-TextCell "D43" "101.25"
$ws.Range("E43").Value = "  +3.09%  "
$ws.Range("D44").Value = "1.822.56"
$ws.Range("E44").Value = "  -1.27%  "
Set-TextCell "D45" "0.00000000113"
$ws.Range("E45").Value = "  +7.22%  "
Set-TextCell "D46" "56.22"
$ws.Range("E46").Value = "  -0.30%  "
Set-TextCell "D47" "1.005"
$ws.Range("E47").Value = "  +0.49%  "
Set-TextCell "D48" "8.103"
$ws.Range("E48").Value = "  +0.94%  "
Set-TextCell "D49" "0.05207"
$ws.Range("E49").Value = "  -0.09%  "
Set-TextCell "D50" "0.4299"
$ws.Range("E50").Value = "  -0.89%  "
$ws.Range("E51").Value = "  +4.17%  "
